$wb = $excel.ActiveWorkbook

# Update "OFF" sheet - Row 3 ("R") values for Short Att, Short Comp, Deep Att, Deep Comp
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 217
$wsOff.Range("C3").Value = 160
$wsOff.Range("D3").Value = 60
$wsOff.Range("E3").Value = 34

# Update "DEF" sheet - Row 3 ("R") values for Short Att, Short Comp, Deep Att, Deep Comp
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 180
$wsDef.Range("C3").Value = 129
$wsDef.Range("D3").Value = 52
$wsDef.Range("E3").Value = 31
